$d = $word.ActiveDocument

$d.Content.Find.Execute("114×3=342", $true, $false, $false, $false, $false, $true, 1, $false, "878×4=3512", 2) | Out-Null
$d.Content.Find.Execute("685×6=4110", $true, $false, $false, $false, $false, $true, 1, $false, "264×9=2376", 2) | Out-Null
$d.Content.Find.Execute("106×9=954", $true, $false, $false, $false, $false, $true, 1, $false, "297×9=2673", 2) | Out-Null
$d.Content.Find.Execute("530×7=3710", $true, $false, $false, $false, $false, $true, 1, $false, "358×4=1432", 2) | Out-Null
$d.Content.Find.Execute("744×6=4464", $true, $false, $false, $false, $false, $true, 1, $false, "380×8=3040", 2) | Out-Null
$d.Content.Find.Execute("661×5=3305", $true, $false, $false, $false, $false, $true, 1, $false, "646×7=4522", 2) | Out-Null
$d.Content.Find.Execute("680×8=5440", $true, $false, $false, $false, $false, $true, 1, $false, "628×4=2512", 2) | Out-Null
$d.Content.Find.Execute("494×5=2470", $true, $false, $false, $false, $false, $true, 1, $false, "839×4=3356", 2) | Out-Null
$d.Content.Find.Execute("127×2=254", $true, $false, $false, $false, $false, $true, 1, $false, "103×4=412", 2) | Out-Null
$d.Content.Find.Execute("315×6=1890", $true, $false, $false, $false, $false, $true, 1, $false, "999×9=8991", 2) | Out-Null
$d.Content.Find.Execute("278×3=834", $true, $false, $false, $false, $false, $true, 1, $false, "342×7=2394", 2) | Out-Null
$d.Content.Find.Execute("800×6=4800", $true, $false, $false, $false, $false, $true, 1, $false, "847×9=7623", 2) | Out-Null
$d.Content.Find.Execute("736×7=5152", $true, $false, $false, $false, $false, $true, 1, $false, "802×5=4010", 2) | Out-Null
$d.Content.Find.Execute("197×2=394", $true, $false, $false, $false, $false, $true, 1, $false, "562×7=3934", 2) | Out-Null
$d.Content.Find.Execute("447×7=3129", $true, $false, $false, $false, $false, $true, 1, $false, "740×4=2960", 2) | Out-Null
$d.Content.Find.Execute("751×3=2253", $true, $false, $false, $false, $false, $true, 1, $false, "345×5=1725", 2) | Out-Null
$d.Content.Find.Execute("799×5=3995", $true, $false, $false, $false, $false, $true, 1, $false, "573×7=4011", 2) | Out-Null
$d.Content.Find.Execute("675×2=1350", $true, $false, $false, $false, $false, $true, 1, $false, "443×7=3101", 2) | Out-Null
$d.Content.Find.Execute("782×8=6256", $true, $false, $false, $false, $false, $true, 1, $false, "208×7=1456", 2) | Out-Null
$d.Content.Find.Execute("455×7=3185", $true, $false, $false, $false, $false, $true, 1, $false, "120×5=600", 2) | Out-Null
$d.Content.Find.Execute("503×8=4024", $true, $false, $false, $false, $false, $true, 1, $false, "718×3=2154", 2) | Out-Null
$d.Content.Find.Execute("436×3=1308", $true, $false, $false, $false, $false, $true, 1, $false, "721×9=6489", 2) | Out-Null
$d.Content.Find.Execute("365×8=2920", $true, $false, $false, $false, $false, $true, 1, $false, "418×3=1254", 2) | Out-Null
$d.Content.Find.Execute("760×8=6080", $true, $false, $false, $false, $false, $true, 1, $false, "546×8=4368", 2) | Out-Null
$d.Content.Find.Execute("899×3=2697", $true, $false, $false, $false, $false, $true, 1, $false, "587×4=2348", 2) | Out-Null
